$d = $word.ActiveDocument

# 1. "Click Alert Box" -> "Type job that already exists "
$d.Content.Find.Execute("Click Alert Box", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Type job that already exists ", 2) | Out-Null

# 2. "Open Alert box" -> "Alert box"
$d.Content.Find.Execute("Open Alert box", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Alert box", 2) | Out-Null

# 3. Remove the existing "_GoBack" bookmark (currently sitting in the
#    "Test Title" header paragraph, right after "buttons and alerts work ").
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# 4. Re-create the "_GoBack" bookmark, collapsed, immediately before the run
#    that now reads "Alert box" in the table cell.
$found = $d.Content.Duplicate
$found.Find.Execute("Alert box", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0) | Out-Null
$bookmarkRange = $d.Range($found.Start, $found.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

Write-Output "edit applied"
